$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: section header "Explicit triggers" - reuse the same formatting
# already used by the other section headers (e.g. A28 "Localizer triggers")
# so no new cell style gets created.
$ws.Range("A28").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Explicit triggers"

$ws.Range("A38").Value = "explicit_start"
$ws.Range("B38").Value = 128

$ws.Range("A39").Value = "explicit_isi"
$ws.Range("B39").Value = 129

$ws.Range("A40").Value = "45_EXP"
$ws.Range("B40").Value = 130

$ws.Range("A41").Value = "45_UEX"
$ws.Range("B41").Value = 131

$ws.Range("A42").Value = "135_EXP"
$ws.Range("B42").Value = 132

$ws.Range("A43").Value = "135_UEX"
$ws.Range("B43").Value = 133

$ws.Range("A44").Value = "100_EXP"
$ws.Range("B44").Value = 134

$ws.Range("A45").Value = "100_UEX"
$ws.Range("B45").Value = 135

$ws.Range("A46").Value = "160_EXP"
$ws.Range("B46").Value = 136

$ws.Range("A47").Value = "160_UEX"
$ws.Range("B47").Value = 137

$ws.Range("A48").Value = "explicit_response"
$ws.Range("B48").Value = 138

$ws.Range("A49").Value = "confidence_response"
$ws.Range("B49").Value = 139

$ws.Range("A37:B50").Select()
$excel.ActiveWindow.ScrollRow = 36
